$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H column values for rows 3,5,7,9,11,13 from -117.901895 to -121.901895
$rows = @(3,5,7,9,11,13)
foreach ($r in $rows) {
    $ws.Range("H$r").Value = -121.901895
}

# Update the view: scroll window back to top-left (A1) and change selection to H13
$ws.Range("H13").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
